# "FEW MORE TEMPLATE UPDATED"
# Update the Policy Holder value on the claims-statement template and
# leave the selection where the user last worked (J3:O3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Policy Holder:- (J2) value cell (J3) MEDGULF TAKAFUL -> ARBIA
$ws.Range("J3").Value = "ARBIA"

# Reflect the cursor/selection left on the sheet after the edit.
[void]$ws.Range("J3:O3").Select()
